$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename two existing "User Interface" sub-requirement rows ---
$ws.Range("A15").Value2 = "Player turn"
$ws.Range("A18").Value2 = "Spin Button"

# --- Rename "Add POINTS to SCORE" sub-requirement to include "(correct answer)" ---
$ws.Range("A44").Value2 = "Add POINTS to SCORE (correct answer)"

# --- Insert 3 new rows right after row 46 (Case 2 of the row-44 sub-requirement), ---
# --- before the old row 47 "Store SCORE for first ROUND during second" header.    ---
$ws.Rows("47:49").Insert()

# Copy formatting from the existing sub-requirement header/case rows as templates,
# then overwrite the text, so style indices (s="1"/"2"/"3") match the rest of the table.
$ws.Range("A50:D50").Copy($ws.Range("A47:D47"))
$ws.Range("B50:D50").Copy($ws.Range("B48:D48"))
$ws.Range("B50:D50").Copy($ws.Range("B49:D49"))

$ws.Range("A47").Value2 = "Subtract POINTS to SCORE (incorrect answer)"
$ws.Range("B48").Value2 = "Case 1"
$ws.Range("B49").Value2 = "Case 2"

# --- Extend Table1 to cover the 3 new rows (was A1:D80, now A1:D83) ---
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:D83"))

# --- Column A width & view changes recorded in the diff ---
$ws.Columns("A").ColumnWidth = 36.5859375

$ws.Application.ActiveWindow.ScrollRow = 13
$ws.Range("B47").Select()
